$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H80").Value = 1017.34784
$ws.Range("I80").Value = 1035.2941
$ws.Range("J80").Value = 966.5
$ws.Range("K80").Value = 3105.8823
$ws.Range("L80").Value = 2899.5
$ws.Range("M80").Value = -2107.8823
$ws.Range("N80").Value = -4895.5
$ws.Range("H83").Value = 1017.34784
$ws.Range("I83").Value = 1035.2941
$ws.Range("J83").Value = 966.5
$ws.Range("K83").Value = 9317.6469
$ws.Range("L83").Value = 8698.5
$ws.Range("M83").Value = -4325.6469
$ws.Range("N83").Value = -18682.5
$ws.Range("H97").Value = 8420
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 761.2857
$ws.Range("I107").Value = 761.2857
$ws.Range("K107").Value = 761.2857
$ws.Range("M107").Value = 1158.7143

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6021.4873
$ws.Range("I32").Value = 3439
$ws.Range("J32").Value = 17827.143
$ws.Range("K32").Value = 3439
$ws.Range("L32").Value = 17827.143
$ws.Range("M32").Value = -3152
$ws.Range("N32").Value = -18401.143
$ws.Range("H45").Value = 4499.75
$ws.Range("I45").Value = 4499.75
$ws.Range("K45").Value = 4499.75
$ws.Range("M45").Value = -4122.75
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H110").Value = 12339.5
$ws.Range("I110").Value = 29999
$ws.Range("J110").Value = 6453
$ws.Range("K110").Value = 29999
$ws.Range("L110").Value = 6453
$ws.Range("M110").Value = -27954
$ws.Range("N110").Value = -10543
$ws.Range("H132").Value = 1132
$ws.Range("I132").Value = 842.6667
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2528.0001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 1.999899999999798
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 1633.8235
$ws.Range("I94").Value = 1281.25
$ws.Range("K94").Value = 1281.25
$ws.Range("M94").Value = -830.25
$ws.Range("H107").Value = 4700
$ws.Range("I107").Value = 4700
$ws.Range("K107").Value = 4700
$ws.Range("M107").Value = -2780
$ws.Range("H134").Value = 3044.9583
$ws.Range("I134").Value = 2185.7
$ws.Range("K134").Value = 6557.099999999999
$ws.Range("M134").Value = -4022.099999999999

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3575.5417
$ws.Range("I31").Value = 3221.2104
$ws.Range("J31").Value = 4922
$ws.Range("K31").Value = 3221.2104
$ws.Range("L31").Value = 4922
$ws.Range("M31").Value = -2926.2104
$ws.Range("N31").Value = -5512
$ws.Range("H34").Value = 3575.5417
$ws.Range("I34").Value = 3221.2104
$ws.Range("J34").Value = 4922
$ws.Range("K34").Value = 3221.2104
$ws.Range("L34").Value = 4922
$ws.Range("M34").Value = -3019.2104
$ws.Range("N34").Value = -5326
$ws.Range("H107").Value = 1180.3
$ws.Range("I107").Value = 752.5
$ws.Range("K107").Value = 752.5
$ws.Range("M107").Value = 1167.5
$ws.Range("H122").Value = 3426.3333
$ws.Range("I122").Value = 3683.3
$ws.Range("K122").Value = 11049.9
$ws.Range("M122").Value = -8599.900000000001
$ws.Range("H132").Value = 1348.4166
$ws.Range("I132").Value = 1355.381
$ws.Range("K132").Value = 4066.143
$ws.Range("M132").Value = -1536.143

$ws = $wb.Worksheets.Item(5)
$ws.Range("H14").Value = 776.75
$ws.Range("I14").Value = 776.75
$ws.Range("K14").Value = 2330.25
$ws.Range("M14").Value = -2157.25
$ws.Range("H92").Value = 370.27777
$ws.Range("I92").Value = 358.18182
$ws.Range("J92").Value = 389.2857
$ws.Range("K92").Value = 1074.54546
$ws.Range("L92").Value = 1167.8571
$ws.Range("M92").Value = 173.45454
$ws.Range("N92").Value = -3663.8571
$ws.Range("H121").Value = 661.3333
$ws.Range("J121").Value = 851.6667
$ws.Range("L121").Value = 2555.0001
$ws.Range("N121").Value = -5175.0001
$ws.Range("H141").Value = 7005.6
$ws.Range("I141").Value = 7005.6
$ws.Range("K141").Value = 21016.8
$ws.Range("M141").Value = -15836.8

$ws = $wb.Worksheets.Item(6)
$ws.Range("H136").Value = 18659
$ws.Range("J136").Value = 18659
$ws.Range("L136").Value = 55977
$ws.Range("N136").Value = -61077

$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 2587.5557
$ws.Range("I40").Value = 2756.6428
$ws.Range("J40").Value = 1995.75
$ws.Range("K40").Value = 2756.6428
$ws.Range("L40").Value = 1995.75
$ws.Range("M40").Value = -2620.6428
$ws.Range("N40").Value = -2267.75
$ws.Range("H68").Value = 2563.625
$ws.Range("I68").Value = 2215.1428
$ws.Range("J68").Value = 5003
$ws.Range("K68").Value = 2215.1428
$ws.Range("L68").Value = 5003
$ws.Range("M68").Value = -1466.1428
$ws.Range("N68").Value = -6501
$ws.Range("H71").Value = 2563.625
$ws.Range("I71").Value = 2215.1428
$ws.Range("J71").Value = 5003
$ws.Range("K71").Value = 11075.714
$ws.Range("L71").Value = 25015
$ws.Range("M71").Value = -7331.714
$ws.Range("N71").Value = -32503
$ws.Range("H122").Value = 4166.6665
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item(8)
$ws.Range("H17").Value = 1750
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -1328
$ws.Range("N17").Value = -2344
$ws.Range("H81").Value = 14777.777
$ws.Range("I81").Value = 19999
$ws.Range("J81").Value = 14125.125
$ws.Range("K81").Value = 39998
$ws.Range("L81").Value = 28250.25
$ws.Range("M81").Value = -38937
$ws.Range("N81").Value = -30372.25
$ws.Range("H84").Value = 14777.777
$ws.Range("I84").Value = 19999
$ws.Range("J84").Value = 14125.125
$ws.Range("K84").Value = 199990
$ws.Range("L84").Value = 141251.25
$ws.Range("M84").Value = -194686
$ws.Range("N84").Value = -151859.25
$ws.Range("H132").Value = 7721.3
$ws.Range("I132").Value = 4037.3333
$ws.Range("J132").Value = 13247.25
$ws.Range("K132").Value = 12111.9999
$ws.Range("L132").Value = 39741.75
$ws.Range("M132").Value = -9581.999899999999
$ws.Range("N132").Value = -44801.75
$ws.Range("H136").Value = 1798.2174
$ws.Range("I136").Value = 1439.9474
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4319.8422
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -1769.8422
$ws.Range("N136").Value = -15600
